# Server: Gacha Table update (Item_List sheet)
#
# - Row 3 (id 100002, Dia): Item_Grade 1 -> 2
# - Row 4 used to be the "기본 상자/Basic_Box" row (id 100011) -> repurposed as
#   a new "마일리지/Money_mileage" currency row (id 100003)
# - Row 5 used to be the "고급상자/High_Box" row (id 100012) -> repurposed as
#   a new "냥냥가챠 뽑기권/Ticket_Basic" row (id 100004)
# - A brand new row is inserted at position 6 for
#   "울트라 냥냥가챠 뽑기권/Ticket_Rare" (id 100005), pushing every row that
#   used to be 6..24 down to 7..25 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Item_List")

# Make room for the new gacha-ticket row: insert a blank row at row 6
# (old rows 6-24 shift down to 7-25, formatting carried from the row below).
$ws.Rows(6).Insert()

# Item_Grade for the Dia currency row changed from 1 to 2.
$ws.Range("E3").Value = 2

# Row 4: "기본 상자" / "Basic_Box" -> "마일리지" / "Money_mileage"
$ws.Range("B4").Value = "마일리지"
$ws.Range("C4").Value = "Money_mileage"

# Row 5: "고급상자" / "High_Box" -> "냥냥가챠 뽑기권" / "Ticket_Basic"
$ws.Range("B5").Value = "냥냥가챠 뽑기권"

# New row 6: "울트라 냥냥가챠 뽑기권" / "Ticket_Rare"
$ws.Range("B6").Value = "울트라 냥냥가챠 뽑기권"

$ws.Range("C5").Value = "Ticket_Basic"
$ws.Range("C6").Value = "Ticket_Rare"

$ws.Range("A4").Value = 100003
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

$ws.Range("A5").Value = 100004
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

$ws.Range("A6").Value = 100005
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45292
$ws.Range("G6").Value = 73050
